$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price / volume data. Two rows (50-51) had their
# coin entirely replaced (Coin name + Link + Price + Volume).
#
# Values that look numeric are written with a leading apostrophe so
# Excel keeps them as literal text (matching the source inlineStr
# cells), instead of silently converting them to numbers.

$ws.Cells.Item(2, 4).Value = '25.437.17'
$ws.Cells.Item(2, 5).Value = '  -1.11%  '

$ws.Cells.Item(3, 4).Value = '1.666.45'
$ws.Cells.Item(3, 5).Value = '  -2.11%  '

$ws.Cells.Item(4, 4).Value = '''0.9984'
$ws.Cells.Item(4, 5).Value = '  -0.88%  '

$ws.Cells.Item(5, 4).Value = '''236.83'
$ws.Cells.Item(5, 5).Value = '  -2.43%  '

$ws.Cells.Item(6, 4).Value = '''0.9993'
$ws.Cells.Item(6, 5).Value = '  -0.60%  '

$ws.Cells.Item(7, 4).Value = '''0.4801'
$ws.Cells.Item(7, 5).Value = '  -3.27%  '

$ws.Cells.Item(8, 4).Value = '''0.2616'
$ws.Cells.Item(8, 5).Value = '  -2.50%  '

$ws.Cells.Item(9, 4).Value = '''0.06152'
$ws.Cells.Item(9, 5).Value = '  +1.26%  '

$ws.Cells.Item(10, 4).Value = '''0.07083'
$ws.Cells.Item(10, 5).Value = '  -1.67%  '

$ws.Cells.Item(11, 4).Value = '1.658.51'
$ws.Cells.Item(11, 5).Value = '  -2.85%  '

$ws.Cells.Item(12, 4).Value = '''14.81'
$ws.Cells.Item(12, 5).Value = '  -0.93%  '

$ws.Cells.Item(13, 4).Value = '''0.5939'
$ws.Cells.Item(13, 5).Value = '  -7.11%  '

$ws.Cells.Item(14, 4).Value = '''4.393'
$ws.Cells.Item(14, 5).Value = '  -6.44%  '

$ws.Cells.Item(15, 4).Value = '''74.50'
$ws.Cells.Item(15, 5).Value = '  -1.38%  '

$ws.Cells.Item(16, 4).Value = '''0.9997'
$ws.Cells.Item(16, 5).Value = '  -0.23%  '

$ws.Cells.Item(17, 4).Value = '''0.9991'
$ws.Cells.Item(17, 5).Value = '  -0.55%  '

$ws.Cells.Item(18, 4).Value = '25.432.67'
$ws.Cells.Item(18, 5).Value = '  -1.34%  '

$ws.Cells.Item(19, 4).Value = '''0.000006752'
$ws.Cells.Item(19, 5).Value = '  +0.31%  '

$ws.Cells.Item(20, 4).Value = '''11.42'
$ws.Cells.Item(20, 5).Value = '  -2.64%  '

$ws.Cells.Item(21, 4).Value = '1.872.93'
$ws.Cells.Item(21, 5).Value = '  -3.97%  '

$ws.Cells.Item(22, 4).Value = '''4.441'
$ws.Cells.Item(22, 5).Value = '  -2.16%  '

$ws.Cells.Item(23, 4).Value = '''8.663'
$ws.Cells.Item(23, 5).Value = '  -0.48%  '

$ws.Cells.Item(24, 4).Value = '''5.344'
$ws.Cells.Item(24, 5).Value = '  -0.69%  '

$ws.Cells.Item(25, 4).Value = '''133.52'
$ws.Cells.Item(25, 5).Value = '  -0.47%  '

$ws.Cells.Item(26, 4).Value = '''15.09'
$ws.Cells.Item(26, 5).Value = '  +0.19%  '

$ws.Cells.Item(27, 4).Value = '''1.398'
$ws.Cells.Item(27, 5).Value = '  -0.36%  '

$ws.Cells.Item(28, 4).Value = '''104.43'
$ws.Cells.Item(28, 5).Value = '  +0.26%  '

$ws.Cells.Item(29, 4).Value = '''1.702'
$ws.Cells.Item(29, 5).Value = '  -2.98%  '

$ws.Cells.Item(30, 4).Value = '''3.994'
$ws.Cells.Item(30, 5).Value = '  +2.38%  '

$ws.Cells.Item(31, 4).Value = '''3.613'
$ws.Cells.Item(31, 5).Value = '  -0.17%  '

$ws.Cells.Item(32, 4).Value = '''0.07652'
$ws.Cells.Item(32, 5).Value = '  -4.98%  '

$ws.Cells.Item(33, 5).Value = '  -6.22%  '

$ws.Cells.Item(34, 4).Value = '''0.9985'
$ws.Cells.Item(34, 5).Value = '  -0.72%  '

$ws.Cells.Item(35, 4).Value = '''2.600'
$ws.Cells.Item(35, 5).Value = '  -2.07%  '

$ws.Cells.Item(36, 4).Value = '''0.6095'
$ws.Cells.Item(36, 5).Value = '  +1.23%  '

$ws.Cells.Item(37, 4).Value = '''0.9456'
$ws.Cells.Item(37, 5).Value = '  -3.87%  '

$ws.Cells.Item(38, 4).Value = '''2.627'
$ws.Cells.Item(38, 5).Value = '  -3.28%  '

$ws.Cells.Item(39, 4).Value = '''0.8581'
$ws.Cells.Item(39, 5).Value = '  -0.32%  '

$ws.Cells.Item(40, 4).Value = '''0.9999'
$ws.Cells.Item(40, 5).Value = '  -0.09%  '

$ws.Cells.Item(41, 4).Value = '''0.01502'
$ws.Cells.Item(41, 5).Value = '  -6.01%  '

$ws.Cells.Item(42, 4).Value = '''1.828'
$ws.Cells.Item(42, 5).Value = '  -5.16%  '

$ws.Cells.Item(43, 4).Value = '''98.58'
$ws.Cells.Item(43, 5).Value = '  -1.91%  '

$ws.Cells.Item(44, 4).Value = '''0.3778'
$ws.Cells.Item(44, 5).Value = '  -1.69%  '

$ws.Cells.Item(45, 5).Value = '  -5.43%  '

$ws.Cells.Item(46, 4).Value = '''0.1122'
$ws.Cells.Item(46, 5).Value = '  -3.96%  '

$ws.Cells.Item(47, 4).Value = '''6.226'
$ws.Cells.Item(47, 5).Value = '  -0.57%  '

$ws.Cells.Item(48, 4).Value = '''0.05253'
$ws.Cells.Item(48, 5).Value = '  -0.20%  '

$ws.Cells.Item(49, 4).Value = '''29.56'
$ws.Cells.Item(49, 5).Value = '  -2.86%  '

$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '''7.352'
$ws.Cells.Item(50, 5).Value = '  -2.60%  '

$ws.Cells.Item(51, 2).Value = 'Decentraland'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(51, 4).Value = '''0.3350'
$ws.Cells.Item(51, 5).Value = '  -2.21%  '
